$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Model formula paragraphs: "LDMC ~" -> "LDMC^(1/3) ~" and a re-wrapped
#    line break before "City_dist" / "Urb_score".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Model: LDMC ~ (1 | Population/Family) + Block + Transect_ID + City_dist +     Transect_ID:City_dist",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "Model: LDMC^(1/3) ~ (1 | Population/Family) + Block + Transect_ID +     City_dist + Transect_ID:City_dist",
    2) | Out-Null

$d.Content.Find.Execute(
    "Model: LDMC ~ (1 | Population/Family) + Block + Transect_ID + Urb_score +     Transect_ID:Urb_score",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "Model: LDMC^(1/3) ~ (1 | Population/Family) + Block + Transect_ID +     Urb_score + Transect_ID:Urb_score",
    2) | Out-Null

# ---------------------------------------------------------------------------
# helper: replace the text of a specific table cell (by 1-based table/row/col)
# ---------------------------------------------------------------------------
function Set-CellText($tableIndex, $rowIndex, $colIndex, $newText) {
    $tbl = $d.Tables.Item($tableIndex)
    $cell = $tbl.Rows.Item($rowIndex).Cells.Item($colIndex)
    $cell.Range.Text = $newText
}

# ---------------------------------------------------------------------------
# 2. Fix the mangled "chi" glyph ("Ï‡" -> "χ") in each table's header row.
#    The header cell holds TWO runs ("chi" + superscript "2"); only replace
#    the mis-encoded first run's text, scoped to that cell, so the
#    superscript "2" run is left untouched.
# ---------------------------------------------------------------------------
$mojibakeA = [char]0x00CF
$mojibakeB = [char]0x2021
$mojibake = "$mojibakeA$mojibakeB"
$chi = [string][char]0x03C7
for ($ti = 1; $ti -le 4; $ti++) {
    $tbl = $d.Tables.Item($ti)
    $cell = $tbl.Rows.Item(1).Cells.Item(3)
    $cell.Range.Find.Execute($mojibake, $false, $false, $false, $false, $false, $true, 1, $false, $chi, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 3. Table 1 (LDMC ~ City_dist, variance-components table)
# ---------------------------------------------------------------------------
Set-CellText 1 2 3 "0.011"
Set-CellText 1 2 5 "0.439"
Set-CellText 1 2 6 "0.459"
Set-CellText 1 3 3 "0.050"
Set-CellText 1 3 5 "0.413"
Set-CellText 1 3 6 "0.4115"
Set-CellText 1 4 4 "0.003"
Set-CellText 1 4 5 "99.148"

# ---------------------------------------------------------------------------
# 4. Table 2 (LDMC ~ City_dist, anova table)
# ---------------------------------------------------------------------------
Set-CellText 2 2 3 "29.173"
Set-CellText 2 2 4 "<0.001***"
Set-CellText 2 3 3 "0.697"
Set-CellText 2 3 4 "0.404"
Set-CellText 2 4 3 "3.266"
Set-CellText 2 4 4 "0.071"
Set-CellText 2 5 3 "0.118"
Set-CellText 2 5 4 "0.731"

# ---------------------------------------------------------------------------
# 5. Table 3 (LDMC ~ Urb_score, variance-components table)
# ---------------------------------------------------------------------------
Set-CellText 3 2 3 "0.013"
Set-CellText 3 2 5 "0.492"
Set-CellText 3 2 6 "0.454"
Set-CellText 3 3 3 "0.211"
Set-CellText 3 3 5 "0.888"
Set-CellText 3 3 6 "0.323"
Set-CellText 3 4 4 "0.003"
Set-CellText 3 4 5 "98.620"

# ---------------------------------------------------------------------------
# 6. Table 4 (LDMC ~ Urb_score, anova table)
# ---------------------------------------------------------------------------
Set-CellText 4 2 3 "28.669"
Set-CellText 4 2 4 "<0.001***"
Set-CellText 4 3 3 "0.687"
Set-CellText 4 3 4 "0.407"
Set-CellText 4 4 3 "0.733"
Set-CellText 4 4 4 "0.392"
Set-CellText 4 5 3 "1.271"
Set-CellText 4 5 4 "0.26"

# ---------------------------------------------------------------------------
# 7. Column widths (gridCol, in twips -> Word COM wants points = twips/20)
#    and header-row heights.
# ---------------------------------------------------------------------------
$d.Tables.Item(1).Columns.Item(3).Width = 961 / 20
$d.Tables.Item(1).Columns.Item(6).Width = 1084 / 20
$d.Tables.Item(1).Rows.Item(1).Height = 615 / 20

$d.Tables.Item(2).Columns.Item(4).Width = 1347 / 20
$d.Tables.Item(2).Rows.Item(1).Height = 571 / 20

$d.Tables.Item(3).Columns.Item(6).Width = 961 / 20
$d.Tables.Item(3).Rows.Item(1).Height = 615 / 20

$d.Tables.Item(4).Columns.Item(4).Width = 1347 / 20
$d.Tables.Item(4).Rows.Item(1).Height = 571 / 20
